# Weekly update: insert two new price rows (week of 2021-11-18) at the top
# of the "Pimiento" data block (row 119), pushing the existing rows down by
# two. This matches the upstream source gaining a new week's observations
# while history is kept sorted with the newest entries first.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at the insertion point; everything at/after row 119
# (previously rows 119:183) shifts down to rows 121:185.
$ws.Rows("119:120").Insert()

# New row 119 - "Zafiro rojo"
$ws.Cells.Item(119,1).Value  = 7
$ws.Cells.Item(119,2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(119,3).Value  = "Ñuble"
$ws.Cells.Item(119,4).Value  = 44518
$ws.Cells.Item(119,5).Value  = 16
$ws.Cells.Item(119,6).Value  = 100112002
$ws.Cells.Item(119,7).Value  = "Pimiento"
$ws.Cells.Item(119,8).Value  = "Zafiro rojo"
$ws.Cells.Item(119,9).Value  = "Primera"
$ws.Cells.Item(119,10).Value = 200
$ws.Cells.Item(119,11).Value = 43000
$ws.Cells.Item(119,12).Value = 44000
$ws.Cells.Item(119,13).Value = 43500
$ws.Cells.Item(119,14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(119,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(119,16).Value = 2900
$ws.Cells.Item(119,17).Value = 15
$ws.Cells.Item(119,18).Value = "Hortaliza"

# New row 120 - "Zafiro verde"
$ws.Cells.Item(120,1).Value  = 7
$ws.Cells.Item(120,2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(120,3).Value  = "Ñuble"
$ws.Cells.Item(120,4).Value  = 44518
$ws.Cells.Item(120,5).Value  = 16
$ws.Cells.Item(120,6).Value  = 100112002
$ws.Cells.Item(120,7).Value  = "Pimiento"
$ws.Cells.Item(120,8).Value  = "Zafiro verde"
$ws.Cells.Item(120,9).Value  = "Primera"
$ws.Cells.Item(120,10).Value = 200
$ws.Cells.Item(120,11).Value = 38000
$ws.Cells.Item(120,12).Value = 39000
$ws.Cells.Item(120,13).Value = 38500
$ws.Cells.Item(120,14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(120,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(120,16).Value = 2567
$ws.Cells.Item(120,17).Value = 15
$ws.Cells.Item(120,18).Value = "Hortaliza"

Write-Host "Dimension now:" $ws.UsedRange.Address()
